# Updated cryptos list on Thu Oct 31 22:40:13 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row 2 - Bitcoin
$ws.Range("D2").Value = "70.363.84"
$ws.Range("E2").Value = "  -2.76%  "

# row 3 - Ethereum
$ws.Range("D3").Value = "2.514.48"
$ws.Range("E3").Value = "  -5.23%  "

# row 4 - TetherUSD
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.17%  "

# row 5 - BNB
$ws.Range("D5").Value = "'575.06"
$ws.Range("E5").Value = "  -3.63%  "

# row 6 - Solana
$ws.Range("D6").Value = "'168.27"
$ws.Range("E6").Value = "  -3.87%  "

# row 7 - USDC
$ws.Range("E7").Value = "  -0.02%  "

# row 8 - XRP
$ws.Range("E8").Value = "  -2.63%  "

# row 9 - LidoStakedEther
$ws.Range("D9").Value = "2.508.31"
$ws.Range("E9").Value = "  -5.50%  "

# row 10 - Dogecoin
$ws.Range("D10").Value = "'0.162"
$ws.Range("E10").Value = "  -4.98%  "

# row 11 - TRON
$ws.Range("E11").Value = "  -0.74%  "

# row 12 - Cardano
$ws.Range("D12").Value = "'0.342"
$ws.Range("E12").Value = "  -3.67%  "

# row 13 - Toncoin
$ws.Range("E13").Value = "  -3.36%  "

# row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "2.969.65"
$ws.Range("E14").Value = "  -5.47%  "

# row 15 - WrappedBTC
$ws.Range("D15").Value = "70.263.76"
$ws.Range("E15").Value = "  -2.78%  "

# row 16 - ShibaInu
$ws.Range("E16").Value = "  -4.02%  "

# row 17 - Avalanche
$ws.Range("D17").Value = "'24.93"
$ws.Range("E17").Value = "  -4.76%  "

# row 18 - WrappedEther
$ws.Range("D18").Value = "2.512.95"
$ws.Range("E18").Value = "  -5.55%  "

# row 19 - Chainlink
$ws.Range("E19").Value = "  -7.67%  "

# row 20 - Uniswap
$ws.Range("D20").Value = "'7.65"
$ws.Range("E20").Value = "  +3.82%  "

# row 21 - BitcoinCash
$ws.Range("D21").Value = "'350.42"
$ws.Range("E21").Value = "  -5.47%  "

# row 22 - Polkadot
$ws.Range("E22").Value = "  -5.83%  "

# row 23 - SuiNetwork
$ws.Range("D23").Value = "'1.97"
$ws.Range("E23").Value = "  -4.24%  "

# row 24 - Dai
$ws.Range("E24").Value = "  +0.03%  "

# row 25 - Litecoin
$ws.Range("D25").Value = "'68.98"
$ws.Range("E25").Value = "  -4.20%  "

# row 26 - NEARProtocol
$ws.Range("D26").Value = "'4.06"
$ws.Range("E26").Value = "  -5.81%  "

# row 27 - Aptos
$ws.Range("D27").Value = "'9.11"
$ws.Range("E27").Value = "  -6.83%  "

# row 28 - WrappedeETH
$ws.Range("D28").Value = "2.639.96"
$ws.Range("E28").Value = "  -5.64%  "

# row 29 - Binance-PegBSC-USD
$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.22%  "

# row 30 - PEPE
$ws.Range("D30").Value = "0.0₃0910"
$ws.Range("E30").Value = "  -6.15%  "

# row 31 - InternetComputer(DFINITY)
$ws.Range("D31").Value = "'7.82"
$ws.Range("E31").Value = "  -3.89%  "

# row 32 - Bittensor
$ws.Range("D32").Value = "'481.72"
$ws.Range("E32").Value = "  -2.61%  "

# row 33 - Fetch.AI
$ws.Range("D33").Value = "'1.29"
$ws.Range("E33").Value = "  -0.11%  "

# row 34 - PancakeSwap
$ws.Range("E34").Value = "  -3.55%  "

# row 35 - FirstDigitalUSD
$ws.Range("D35").Value = "'0.998"
$ws.Range("E35").Value = "  -0.17%  "

# row 36 - Monero
$ws.Range("D36").Value = "'157.33"
$ws.Range("E36").Value = "  -2.97%  "

# row 37 - Kaspa
$ws.Range("E37").Value = "  +3.36%  "

# row 38 - WhiteBITCoin
$ws.Range("D38").Value = "'18.90"
$ws.Range("E38").Value = "  -0.20%  "

# row 39 - EthereumClassic
$ws.Range("D39").Value = "'18.57"
$ws.Range("E39").Value = "  -4.65%  "

# row 40 - USDe
$ws.Range("E40").Value = "  -0.02%  "

# row 41 - PolygonEcosystemToken
$ws.Range("D41").Value = "'0.320"
$ws.Range("E41").Value = "  -3.58%  "

# row 42 - Stacks->RenderToken
$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'4.73"
$ws.Range("E42").Value = "  -5.28%  "

# row 43 - RenderToken->Stacks
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").Value = "'1.64"
$ws.Range("E43").Value = "  -6.88%  "

# row 44 - ImmutableX
$ws.Range("E44").Value = "  -13.35%  "

# row 45 - dogwifhat
$ws.Range("D45").Value = "'2.38"
$ws.Range("E45").Value = "  -7.94%  "

# row 46 - OKB
$ws.Range("D46").Value = "'38.28"
$ws.Range("E46").Value = "  -2.44%  "

# row 47 - Aave
$ws.Range("D47").Value = "'142.82"
$ws.Range("E47").Value = "  -8.11%  "

# row 48 - Filecoin
$ws.Range("E48").Value = "  -5.75%  "

# row 49 - ARBITRUM
$ws.Range("E49").Value = "  -5.39%  "

# row 50 - Optimism
$ws.Range("E50").Value = "  -6.21%  "

# row 51 - Mantle
$ws.Range("E51").Value = "  -1.43%  "
